$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "245.43"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "27.17"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "3.38%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.105"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.41%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05695"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.63%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.514"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.49%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8196"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.83%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8593"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "1.92%"
$ws.Range("E8").Style = "Normal"
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1335"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.28%"
$ws.Range("E9").Style = "Normal"
$ws.Range("B10").Value = "MandalaExchangeToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06954"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-0.68%"
$ws.Range("E10").Style = "Normal"
$ws.Range("B11").Value = "BitrueCoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.02867"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.49%"
$ws.Range("E11").Style = "Normal"
$ws.Range("B12").Value = "BitMartToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09393"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.05%"
$ws.Range("E12").Style = "Normal"
$ws.Range("B13").Value = "BitForexToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.001518"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.62%"
$ws.Range("E13").Style = "Normal"
$ws.Range("B14").Value = "CoinExToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.04051"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-13.03%"
$ws.Range("E14").Style = "Normal"
$ws.Range("B15").Value = "One"
$ws.Range("C15").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006019"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-93.99%"
$ws.Range("E15").Style = "Normal"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006214"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.98%"
$ws.Range("E16").Style = "Normal"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.509"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-2.69%"
$ws.Range("E17").Style = "Normal"
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.009"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.11%"
$ws.Range("E18").Style = "Normal"
$ws.Range("B19").Value = "BTSEToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.317"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "12.74%"
$ws.Range("E19").Style = "Normal"
$ws.Range("B20").Value = "BitpandaEcosystemToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3165"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "1.23%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.03213"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.29%"
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-1.81%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.563"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-5.15%"
$ws.Range("E23").Style = "Normal"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-2.36%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004477"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-2.37%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.00009896"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "3.11%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0001449"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "-25.23%"
$ws.Range("E28").Style = "Normal"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.78%"
$ws.Range("E40").Style = "Normal"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.005938"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-3.74%"
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1057"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-0.26%"
$ws.Range("E42").Style = "Normal"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002378"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-4.87%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.009720"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "17.57%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005159"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-4.23%"
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.02%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-8.18%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002508"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-3.30%"
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.02%"
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.02%"
$ws.Range("E50").Style = "Normal"
